# Applies the "Actualización automática" update:
# RIOS CARRION ANGEL BENIGNO / FERRIACABADOS MACONSE now has a PORCELANATO sale
# of 7529.26 recorded for julio, which ripples into the monthly totals and the
# compliance ("CUMPLIMIENTO MENSUAL") summary sheet.

$wb = $excel.ActiveWorkbook

$wsVentasGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual  = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ---
# Row 13 = FERRIACABADOS MACONSE, column M = PORCELANATO
$wsVentasGrupo.Range("M13").Value = 7529.26
# Row 24 = totals row; column M counter of non-zero entries (2 -> 3 of 22)
$wsVentasGrupo.Range("M24").Value = "3 de 22"

# --- VENTA MENSUAL ---
# Row 13 = FERRIACABADOS MACONSE, column F = julio
$wsVentaMensual.Range("F13").Value = 7529.26
# Row 24 = totals row, column F = julio total
$wsVentaMensual.Range("F24").Value = 17113.46

# --- CUMPLIMIENTO MENSUAL ---
# Widen column D ("VENTA") slightly to fit the new figure.
# (The runtime's ColumnWidth -> stored OOXML width conversion adds an offset
#  of 0.8333..; subtract it so the saved <col width="14".../> matches exactly.)
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.166666666666666

# Row 16 = PORCELANATO group
$wsCumplimiento.Range("D16").Value = 17113.46
$wsCumplimiento.Range("E16").Value = 21643.08
$wsCumplimiento.Range("F16").Value = 0.4415631529543143

# Row 19 = TOTAL row
$wsCumplimiento.Range("D19").Value = 17113.46
$wsCumplimiento.Range("E19").Value = 41109.54386304603
$wsCumplimiento.Range("F19").Value = 0.2939295272407245
